# Apply the "added new s tables" edit:
#  - insert new "Table S6" and "Table S7" descriptive paragraphs (each
#    followed by a blank spacer paragraph) right before the existing
#    "Table S10" paragraph
#  - reword the "Table S10" / "Table S11" paragraphs: capitalize the
#    leading "change" -> "Change" and append "; linear mixed-effects
#    model." in place of the old trailing period.

$d = $word.ActiveDocument
$d.TrackRevisions = $false

# ---------------------------------------------------------------------
# 1. Insert the Table S6 (CuSO4) and Table S7 (NaCl) paragraphs before
#    the paragraph that currently begins "Table S10:".
# ---------------------------------------------------------------------

# Locate the "Table S10" paragraph by scanning (robust to any earlier
# paragraph-count drift).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Table S10:")) {
        $targetIndex = $i
        break
    }
}

$insertPos = $d.Paragraphs.Item($targetIndex).Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

# Build the two new paragraphs + their trailing blank spacer paragraphs
# as one text blob (Word paragraph marks = "`r"), then go back and apply
# character formatting (bold label, subscript "4").
$s6 = "Table S6: Cost of generalization for treatments EH0_40, EH20_60, EH0_80, EH40_80 relative to constant environment treatments EH0 and EH80 in 0%, 80% chemical stress for the CuSO4 dataset; linear model`r`r"
$s7 = "Table S7: Cost of generalization for treatments EH0_40, EH20_60, EH0_80, EH40_80 relative to constant environment treatments EH0 and EH80 in 0%, 80% chemical stress for the NaCl dataset; linear model`r`r"

$insertRange.InsertBefore($s6 + $s7)

# Re-find the "Table S10" paragraph again (index shifted by the 4 new
# paragraphs: Table S6 text, blank, Table S7 text, blank).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Table S10:")) {
        $targetIndex = $i
        break
    }
}

$s6Index = $targetIndex - 4
$s7Index = $targetIndex - 2

# --- bold the "Table S6:" label ---
$s6Para = $d.Paragraphs.Item($s6Index)
$s6Start = $s6Para.Range.Start
$boldRange = $d.Range($s6Start, $s6Start + 9)
$boldRange.Font.Bold = $true

# --- subscript the "4" in "CuSO4" for Table S6 ---
$s6Para = $d.Paragraphs.Item($s6Index)
$s6Text = $s6Para.Range.Text
$cusoOffset = $s6Text.IndexOf("CuSO4")
$subStart = $s6Para.Range.Start + $cusoOffset + 4
$subRange = $d.Range($subStart, $subStart + 1)
$subRange.Font.Subscript = $true

# --- bold the "Table S7:" label ---
$s7Para = $d.Paragraphs.Item($s7Index)
$s7Start = $s7Para.Range.Start
$boldRange2 = $d.Range($s7Start, $s7Start + 9)
$boldRange2.Font.Bold = $true

# ---------------------------------------------------------------------
# 2. Reword the (now shifted) "Table S10" / "Table S11" paragraphs.
# ---------------------------------------------------------------------

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Table S10:") -or $t.StartsWith("Table S11:")) {
        $r = $p.Range
        $r.Find.Execute(" change in variance fitness", $true, $false, $false, $false, $false, $true, 1, $false, " Change in variance fitness", 1)
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Table S10:") -and $t.Contains("NaCl dataset.")) {
        $r = $p.Range
        $r.Find.Execute("NaCl dataset.", $true, $false, $false, $false, $false, $true, 1, $false, "NaCl dataset; linear mixed-effects model.", 1)
    }
    if ($t.StartsWith("Table S11:") -and $t.Contains("CuSO4 dataset.")) {
        $r = $p.Range
        $r.Find.Execute("CuSO4 dataset.", $true, $false, $false, $false, $false, $true, 1, $false, "CuSO4 dataset; linear mixed-effects model.", 1)
    }
}

Write-Output "done"
